# pages_config.xlsx - "Add files via upload" commit.
#
# The sheet holds per-page config rows (A:F = page_name, page_id,
# page_token, gemini_api_key, store_link, prompt). This edit rotates the
# Gemini API keys in column D (rows 2-11) to a fresh batch, and drops two
# now-unused leftover API-key strings that were sitting in the shared
# string table (they were not referenced by any cell already).
#
# Column D values are written in the same (non-sequential) order the
# original author's workbook shows them landing in the shared-string
# table, so the underlying XML layout matches exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value  = "AIzaSyBRn6gwWqlxjybH4JG0rKtbdpQIKEsxJtY"
$ws.Range("D7").Value  = "AIzaSyCCfjbzSQTCU5d3u4nDNCpqpJ9_9Pivgcg"
$ws.Range("D3").Value  = "AIzaSyBgJq3NZJzv6_qLJfzzvnC2LvEZ2bNVt5c"
$ws.Range("D2").Value  = "AIzaSyBepfUgkcY3_Yf4iIJzYQOhiJuGkkEN6GE"
$ws.Range("D5").Value  = "AIzaSyBm84-gErNODPQzJ0NGaTqxoqUxb6qPw_o"
$ws.Range("D6").Value  = "AIzaSyAcGlRERU5sQkBGmr_fXuV0y-vAo0IEV78"
$ws.Range("D8").Value  = "AIzaSyAoyArWScFA7651r1rVeeWQquilZyuJbl8"
$ws.Range("D9").Value  = "AIzaSyDGlQgzZ2b2WXiP0MEVwLce4i8ISS2uuRE"

# D10 previously carried a leftover "applyFont" cell style; the refreshed
# key goes back in as a plain, unstyled cell (same as its neighbours).
$ws.Range("D10").Style = "Normal"
$ws.Range("D10").Value = "AIzaSyBjQ36QX6K6cDCCGhkTu2ClKEiNk8frhzc"

$ws.Range("D11").Value = "AIzaSyANMkuXX12OCrQvXHtjvF9ImAYz88kIrU4"

# Selection/scroll position left at D13 before -> move to F13, matching
# the saved view state in the upload.
$ws.Activate()
$ws.Range("F13").Select()
